# Correcciones en el dashboard
# Update "Salario básico" values (column C) in Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C4"  = 200
    "C6"  = 218
    "C8"  = 240
    "C10" = 264
    "C12" = 292
    "C14" = 318
    "C16" = 340
    "C18" = 354
    "C20" = 366
    "C22" = 375
    "C24" = 386
    "C26" = 394
    "C28" = 400
    "C32" = 425
    "C34" = 450
    "C36" = 460
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
